$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (2-9) down by one to make room for the new
# weekly record at row 2. Copy whole rows bottom-up so nothing gets
# overwritten before it is moved, and so each row keeps its own
# formatting (e.g. the date-number-format style on column D).
for ($r = 9; $r -ge 2; $r--) {
    $srcRow = $r
    $dstRow = $r + 1
    $src = $ws.Range("A" + $srcRow + ":R" + $srcRow)
    $dst = $ws.Range("A" + $dstRow + ":R" + $dstRow)
    $src.Copy($dst)
}

# Write the new record into row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").Value = 45282
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 100112017
$ws.Range("G2").Value = "Corazón de apio"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Segunda"
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 900
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = 940
$ws.Range("N2").Value = "$/docena de matas"
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 157
$ws.Range("Q2").Value = 6
$ws.Range("R2").Value = "Hortaliza"
